# Update cryptocurrency price/volume data per the Mon Apr 15 23:25:04 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.390.61'
$ws.Range('E2').Value = '  -3.23%  '

$ws.Range('D3').Value = '3.110.50'
$ws.Range('E3').Value = '  -1.63%  '

$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').Value = "'552.69"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.77%  '

$ws.Range('D6').Value = "'138.81"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.61%  '

$ws.Range('E7').Value = '  +0.24%  '

$ws.Range('D8').Value = '3.104.76'
$ws.Range('E8').Value = '  -1.88%  '

$ws.Range('D9').Value = "'0.496"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.92%  '

$ws.Range('D10').Value = "'0.162"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.27%  '

$ws.Range('D11').Value = "'6.59"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.72%  '

$ws.Range('D12').Value = "'0.458"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.79%  '

$ws.Range('D13').Value = "'35.23"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.72%  '

$ws.Range('D14').Value = "'0.0000219"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.30%  '

$ws.Range('D15').Value = '3.613.96'
$ws.Range('E15').Value = '  -1.52%  '

$ws.Range('D16').Value = '63.448.70'
$ws.Range('E16').Value = '  -2.91%  '

$ws.Range('D17').Value = "'0.111"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.04%  '

$ws.Range('D18').Value = '3.118.79'
$ws.Range('E18').Value = '  -1.35%  '

$ws.Range('D19').Value = "'507.46"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.39%  '

$ws.Range('D20').Value = "'6.68"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.28%  '

$ws.Range('D21').Value = "'13.58"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.41%  '

$ws.Range('D22').Value = "'0.708"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.96%  '

$ws.Range('D23').Value = "'7.26"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.01%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = "'78.18"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.42%  '

$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = "'12.41"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.81%  '

$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('D27').Value = "'2.76"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.06%  '

$ws.Range('D28').Value = "'8.30"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.20%  '

$ws.Range('D29').Value = "'1.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.22%  '

$ws.Range('D30').Value = "'1.95"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -10.61%  '

$ws.Range('D31').Value = "'26.36"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.62%  '

$ws.Range('E32').Value = '  -8.03%  '

$ws.Range('E33').Value = '  -2.63%  '

$ws.Range('D34').Value = "'58.34"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +8.83%  '

$ws.Range('D35').Value = "'530.37"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -10.50%  '

$ws.Range('D36').Value = "'6.00"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.50%  '

$ws.Range('D37').Value = "'5.24"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.71%  '

$ws.Range('D38').Value = "'0.0415"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.33%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '3.077.81'
$ws.Range('E39').Value = '  -0.16%  '

$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = "'0.0795"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.33%  '

$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = "'0.121"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.69%  '

$ws.Range('D42').Value = "'2.75"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -9.64%  '

$ws.Range('D43').Value = "'8.15"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.29%  '

$ws.Range('D44').Value = "'0.254"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.18%  '

$ws.Range('E45').Value = '  +74.42%  '

$ws.Range('E46').Value = '  +0.08%  '

$ws.Range('D47').Value = "'2.06"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.25%  '

$ws.Range('D48').Value = "'123.07"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.26%  '

$ws.Range('D49').Value = "'24.42"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.27%  '

$ws.Range('D50').Value = "'0.107"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.16%  '

$ws.Range('D51').Value = '0.0₃0510'
$ws.Range('E51').Value = '  -8.13%  '
